$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, avoiding Excel's automatic
# numeric coercion for strings that look like numbers (e.g. '528.50'),
# while leaving the cell's style/number-format untouched afterwards.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '59.335.98'
$ws.Range('E2').Value = '  -0.33%  '

$ws.Range('D3').Value = '2.645.00'
$ws.Range('E3').Value = '  -0.42%  '

$ws.Range('E4').Value = '  -0.23%  '

Set-TextValue 'D5' '528.50'
$ws.Range('E5').Value = '  +1.79%  '

Set-TextValue 'D6' '145.23'
$ws.Range('E6').Value = '  -1.34%  '

Set-TextValue 'D7' '0.997'
$ws.Range('E7').Value = '  -0.28%  '

Set-TextValue 'D8' '0.570'
$ws.Range('E8').Value = '  +0.10%  '

Set-TextValue 'D9' '6.68'
$ws.Range('E9').Value = '  -2.93%  '

$ws.Range('E10').Value = '  +1.60%  '

Set-TextValue 'D11' '0.337'
$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('E12').Value = '  +0.98%  '

$ws.Range('D13').Value = '3.111.85'
$ws.Range('E13').Value = '  -0.46%  '

$ws.Range('D14').Value = '59.282.29'
$ws.Range('E14').Value = '  -0.46%  '

$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D15' '21.10'
$ws.Range('E15').Value = '  +0.24%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D16' '0.0000137'
$ws.Range('E16').Value = '  +1.09%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.671.64'
$ws.Range('E17').Value = '  +1.15%  '

Set-TextValue 'D18' '342.46'
$ws.Range('E18').Value = '  +0.90%  '

Set-TextValue 'D19' '4.46'
$ws.Range('E19').Value = '  +0.65%  '

Set-TextValue 'D20' '10.60'
$ws.Range('E20').Value = '  +2.60%  '

Set-TextValue 'D21' '6.36'
$ws.Range('E21').Value = '  +0.95%  '

$ws.Range('E22').Value = '  +0.04%  '

Set-TextValue 'D23' '65.30'
$ws.Range('E23').Value = '  +3.45%  '

Set-TextValue 'D24' '0.419'
$ws.Range('E24').Value = '  +1.80%  '

Set-TextValue 'D25' '0.168'
$ws.Range('E25').Value = '  +0.03%  '

$ws.Range('E26').Value = '  -0.68%  '

Set-TextValue 'D27' '7.28'
$ws.Range('E27').Value = '  +1.75%  '

$ws.Range('D28').Value = '0.0₃0803'
$ws.Range('E28').Value = '  +0.18%  '

$ws.Range('E29').Value = '  -3.85%  '

$ws.Range('E30').Value = '  -0.16%  '

Set-TextValue 'D31' '1.62'
$ws.Range('E31').Value = '  +2.39%  '

Set-TextValue 'D32' '18.98'
$ws.Range('E32').Value = '  +1.10%  '

Set-TextValue 'D33' '150.27'
$ws.Range('E33').Value = '  +0.27%  '

Set-TextValue 'D34' '4.23'
$ws.Range('E34').Value = '  +1.64%  '

Set-TextValue 'D35' '1.21'
$ws.Range('E35').Value = '  +0.81%  '

Set-TextValue 'D36' '0.919'
$ws.Range('E36').Value = '  +2.05%  '

Set-TextValue 'D37' '0.875'
$ws.Range('E37').Value = '  -0.26%  '

Set-TextValue 'D38' '1.51'
$ws.Range('E38').Value = '  +0.86%  '

Set-TextValue 'D39' '36.61'
$ws.Range('E39').Value = '  -0.73%  '

Set-TextValue 'D40' '3.67'
$ws.Range('E40').Value = '  +2.50%  '

Set-TextValue 'D41' '0.996'
$ws.Range('E41').Value = '  -0.42%  '

$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D42' '274.06'
$ws.Range('E42').Value = '  -0.54%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D43' '0.605'
$ws.Range('E43').Value = '  -4.24%  '

$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D44' '0.0976'
$ws.Range('E44').Value = '  +0.18%  '

Set-TextValue 'D45' '19.44'
$ws.Range('E45').Value = '  -2.02%  '

Set-TextValue 'D46' '0.0539'
$ws.Range('E46').Value = '  +0.78%  '

$ws.Range('E47').Value = '  +1.16%  '

$ws.Range('D48').Value = '2.052.67'
$ws.Range('E48').Value = '  -0.47%  '

Set-TextValue 'D49' '4.83'
$ws.Range('E49').Value = '  +0.75%  '

$ws.Range('E50').Value = '  +0.57%  '

Set-TextValue 'D51' '18.99'
$ws.Range('E51').Value = '  -0.09%  '
